$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9 updates
$ws.Range("D9").Value = "국내 대학 교육의 비참한 실태 (2)"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/korean-higher-education-miserable-2/#utm_source=rss&utm_medium=rss&utm_campaign=korean-higher-education-miserable-2"

# Row 26 update
$ws.Range("D26").Value = "ai plus(est soft)"

# Row 51 updates
$ws.Range("D51").Value = "[vim] 윈도우 PC에 vim 설치하기"
$ws.Range("E51").Value = "https://bskyvision.com/1301"
